$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.699881911277771
$ws.Range("B1").Value = 2.89067530632019
$ws.Range("C1").Value = 6.075221538543701
$ws.Range("D1").Value = 2.207178592681885
$ws.Range("E1").Value = 0.7949870228767395
